$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Marca" header column in N1 (extends the table from M to N)
$ws.Range("N1").Value = "Marca"

# Update the active selection to N1 (matches the authored selection change)
$ws.Range("N1").Select()
